$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.988.25"
$ws.Range("D3").Value = "3.345.04"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'606.31"
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("D6").Value = "'143.54"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "3.344.61"
$ws.Range("E8").Value = "  +2.61%  "
$ws.Range("D9").Value = "'0.521"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").Value = "'0.151"
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("D11").Value = "'5.59"
$ws.Range("E11").Value = "  +4.75%  "
$ws.Range("D12").Value = "'0.472"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").Value = "'0.0000250"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").Value = "'35.33"
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("D15").Value = "3.891.74"
$ws.Range("E15").Value = "  +2.69%  "
$ws.Range("D16").Value = "'0.121"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "3.324.15"
$ws.Range("E17").Value = "  +2.26%  "
$ws.Range("D18").Value = "64.061.43"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").Value = "'6.90"
$ws.Range("E19").Value = "  +1.58%  "
$ws.Range("D20").Value = "'483.38"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").Value = "'14.21"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "'0.741"
$ws.Range("E22").Value = "  +1.77%  "
$ws.Range("D23").Value = "'8.05"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("D24").Value = "'13.97"
$ws.Range("E24").Value = "  +5.09%  "
$ws.Range("D25").Value = "'85.08"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").Value = "'2.79"
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'8.33"
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "'7.24"
$ws.Range("E30").Value = "  -3.50%  "
$ws.Range("D31").Value = "'2.17"
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("D32").Value = "'29.00"
$ws.Range("E32").Value = "  +4.56%  "
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("D34").Value = "'2.55"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").Value = "'1.11"
$ws.Range("E35").Value = "  +1.01%  "
$ws.Range("D36").Value = "'6.11"
$ws.Range("E36").Value = "  +3.23%  "
$ws.Range("D37").Value = "0.0₃0758"
$ws.Range("E37").Value = "  +5.61%  "
$ws.Range("D38").Value = "'52.52"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D39").Value = "'0.0401"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("D40").Value = "3.138.05"
$ws.Range("E40").Value = "  +4.79%  "
$ws.Range("D41").Value = "'433.81"
$ws.Range("E41").Value = "  +2.59%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.119"
$ws.Range("E42").Value = "  +6.30%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.78"
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D44").Value = "'8.39"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").Value = "'0.269"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").Value = "'2.26"
$ws.Range("E46").Value = "  +4.19%  "
$ws.Range("D47").Value = "'26.61"
$ws.Range("E47").Value = "  +2.34%  "
$ws.Range("D48").Value = "'36.55"
$ws.Range("E48").Value = "  +7.80%  "
$ws.Range("D50").Value = "'2.33"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("E51").Value = "  -0.50%  "
